# Apply highlight colors to the two Draft sections per the commit:
# "Made desions regarding the 2 drafts via highlights"
#
# WdColorIndex values used (confirmed against this runtime's OOXML mapping):
#   6 = wdRed    -> w:val="red"
#   4 = wdBrightGreen -> w:val="green"
#   7 = wdYellow -> w:val="yellow"
#
# For plain paragraphs (no paragraph-mark formatting change in the diff) we set
# HighlightColorIndex directly on the paragraph Range, which only stamps the
# runs' rPr.
# For paragraphs where the diff also stamps a <w:pPr><w:rPr><w:highlight/></w:rPr>
# (the bulleted/numbered list items, and the "application is deployed" paragraph)
# we go through the Range's Font object, which additionally stamps the
# paragraph-mark run properties inside <w:pPr>.

$d = $word.ActiveDocument

$wdRed    = 6
$wdGreen  = 4
$wdYellow = 7

function Set-ParaHighlight($paraIndex, $color) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.HighlightColorIndex = $color
}

function Set-ParaHighlightWithMark($paraIndex, $color) {
    $p = $d.Paragraphs($paraIndex)
    $p.Range.Font.HighlightColorIndex = $color
}

# ---------- Draft 1: Concise IEEE Style ----------
Set-ParaHighlight 3 $wdRed      # The current version of the system is purely front-end...
Set-ParaHighlight 4 $wdGreen    # The interface utilizes the following front-end libraries and frameworks:

Set-ParaHighlightWithMark 5 $wdYellow   # Bootstrap 5 RTL: for responsive layout and component styling
Set-ParaHighlightWithMark 6 $wdYellow   # Font Awesome: for iconography
Set-ParaHighlightWithMark 7 $wdYellow   # HTML5/CSS3: for structure and presentation

Set-ParaHighlight 8 $wdRed      # The system is hosted on GitHub Pages, which serves static content...

# ---------- Draft 2: Academic Descriptive Style ----------
Set-ParaHighlight 12 $wdRed     # This system is implemented entirely as a front-end web application...
Set-ParaHighlight 13 $wdYellow  # It interfaces with the following software libraries and platforms:

Set-ParaHighlightWithMark 14 $wdGreen   # Bootstrap 5 RTL: Provides responsive UI components...
Set-ParaHighlightWithMark 15 $wdGreen   # Font Awesome: Enables the use of scalable icons...
Set-ParaHighlightWithMark 16 $wdGreen   # HTML5 and CSS3: Core web technologies for layout...

Set-ParaHighlightWithMark 17 $wdRed     # The application is deployed using GitHub Pages...

Set-ParaHighlight 18 $wdGreen   # Future versions may integrate backend services and databases...
